# Insert a new weekly price record at row 41 for
# "Terminal Hortofrutícola Agro Chillán" / Chirimoya, pushing the
# existing rows 41-46 down to 42-47.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 41:46 down to 42:47 by inserting a new row above row 41.
$ws.Rows.Item(41).Insert()

# Excel's Insert() usually copies formatting from the row above; make sure
# the date cell (column D) keeps the same date-number style used by the
# other rows in this block (style index 2 -> numFmtId 165) by copying it
# from the row that used to be directly above (row 40, still at 40).
$ws.Range("D40").Copy()
$ws.Range("D41").PasteSpecial(-4122) # xlPasteFormats

# Populate the new row with the values from the diff.
$ws.Range("A41").Value = 7
$ws.Range("B41").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C41").Value = "Ñuble"
$ws.Range("D41").Value = 45275
$ws.Range("E41").Value = 16
$ws.Range("F41").Value = "Fruta"
$ws.Range("G41").Value = 100107
$ws.Range("H41").Value = "Otros"
$ws.Range("I41").Value = 100107002
$ws.Range("J41").Value = "Chirimoya"
$ws.Range("K41").Value = "Cultivar IV Región"
$ws.Range("L41").Value = "Primera"
$ws.Range("M41").Value = 30
$ws.Range("N41").Value = 19000
$ws.Range("O41").Value = 19000
$ws.Range("P41").Value = 19000
$ws.Range("Q41").Value = "$/bandeja 10 kilos"
$ws.Range("R41").Value = "Provincia de Limarí"
$ws.Range("S41").Value = 1900
$ws.Range("T41").Value = 10
